# edit.ps1 -- apply the commit's changes via PowerPoint COM-interop
#
# 1) Bump the cached "datetimeFigureOut" field text from 17.03.2018 to
#    24.03.2018 everywhere it is stored (the slide master and every
#    slide layout each carry their own cached copy of the field).
# 2) On slide 2 ("Argon use scoreboard"), split the
#    "Out of Order Execution   (E)" run into two runs so the trailing
#    "(E)" becomes "(X)".
# 3) On the same slide, coalesce the four runs that spell out
#    "WPIPE stores WB bit and register that we are going to write"
#    into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached date field text (master + all custom layouts).
# ---------------------------------------------------------------------
function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "17.03.2018") {
                $tr.Text = "24.03.2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) & 3) Text tweaks on slide 2's content placeholder.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(2)
$body = $slide.Shapes.Item(2)
$bodyRange = $body.TextFrame.TextRange

# --- (2) "Out of Order Execution   (E)" -> "...   " + "(X)" -----------
$para5 = $bodyRange.Paragraphs(5, 1)
$tail = $para5.Characters(26, 3)   # the "(E)" substring
$tail.Text = "(X)"

# --- (3) merge the WPIPE-stores run sequence into a single run --------
$para10 = $bodyRange.Paragraphs(10, 1)
# The target text is identical to the already-concatenated paragraph
# text, so assign a throwaway value first to force the run rebuild,
# then set the real text -- otherwise a same-value assignment is a
# no-op and the four original runs would be left untouched.
$para10.Text = "~~temp~~"
$para10b = $bodyRange.Paragraphs(10, 1)
$para10b.Text = "WPIPE stores WB bit and register that we are going to write"
